$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = "[name=""Greatmouth Mob""]  Wait—as I was saying that, it looks like 'Hoarhair' has taken a point off 'Springwater'—oh! And 'Drifter' Dorka barges in to pick them off! Now that’s what I call productive!   `n"
$ws.Range("C7").Value = "[name=""'Gałązka' Knight""]  That’s all you’re capable of? How the hell did you beat Ingra?! `n"
$ws.Range("C9").Value = "[name=""'Gałązka' Knight""]  Heh, I’ll tear you right off that scoreboard, 'Nearl!' Watch this!`n"
$ws.Range("C11").Value = "[name=""'Gałązka' Knight""]  Uagh—?! `n"
$ws.Range("C14").Value = "[name=""Greatmouth Mob""]  And right after 'Gałązka' comes 'Flametail'! Whoa, way-way-way, how many times has someone come for the Nascent Nearl now? You gotta get your own show, Maria!  `n"
$ws.Range("C35").Value = "[name=""'Gałązka' Knight""]  Urrrrraaahhh! You stinking little squirrel—! `n"
$ws.Range("C37").Value = "[name=""'Gałązka' Knight""]  Hah! Knocked away your weapon! What are you gonna do now, huh?! `n"
$ws.Range("C39").Value = "[name=""'Gałązka' Knight""]  Don’t waste your breath! You’ll meet your end right here! `n"
$ws.Range("C41").Value = "[name=""'Gałązka' Knight""]  What the hell? Arts—no, you can’t be, you’re not even armed!`n"
$ws.Range("C42").Value = "[name=""'Gałązka' Knight""]  No, you... you... you can’t be—`n"
$ws.Range("C54").Value = "[name=""Passing Fan""]  When the Blood Knight took the title, they went 'yeah, we’ll officially let Infected be knights'—but when the Radiant Knight got infected, they chased her completely out of Kazimierz. Doesn’t that suck?`n"
$ws.Range("C59").Value = "[name=""Passing Fan""]  Welcome to the bottom line. They weren’t gonna give the Infected any respect for a single second. 'Just throw them into the machine. Make them one more symbol of our progress.'`n"
$ws.Range("C64").Value = "[name=""Passing Fan""]  Let’s call the points I’ve already taken a fee, and here’s a lesson on behalf of 'Whislash'. From the match where you beat Ingra, I can tell you still don’t understand what the actual significance of the Major is.`n"
$ws.Range("C86").Value = "Lightly, like she’s just messing with me, 'Flametail' taps against my shield.`n"
$ws.Range("C87").Value = "[name=""'Flametail' Knight""]  Maybe you’d have taken one more point off me here. Maybe not. Who can say, Maria?`n"
$ws.Range("C89").Value = "[name=""'Flametail' Knight""]  Hahahaha—you know, from the stands, it’s hard to tell just how fiercely stubborn you are, Nearlie.`n"
$ws.Range("C92").Value = "[name=""Greatmouth Mob""]  Fourth place! Earning ten points, the arena was his to walk, and he avoided almost any direct battle! Our survivor with smarts, 'Gałązka' Daniel!  `n"
$ws.Range("C94").Value = "[name=""Greatmouth Mob""]  Second place! In the starting skirmishes, he beat three knights with his own hands! He made the arena’s corner his own with his tremendous shields! The walking fortress, 'Limestone' Marko!   `n"
$ws.Range("C95").Value = "[name=""Greatmouth Mob""]  And our champion tonight, with twenty-two fresh points to her name! A near uninterrupted combo of victories! The super-rookie who never drags her feet—'Flametail' Sonna! `n"
